$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$ws.Range("B2").Value = -0.09914768942982544
$ws.Range("C2").Value = 0.6760636891099804
$ws.Range("D2").Value = 1.058759192261776
$ws.Range("E2").Value = 1.02896024814459
$ws.Range("F2").Value = 1.052237025233689
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = -0.1143386297721907
$ws.Range("C3").Value = 0.6889840097079422
$ws.Range("D3").Value = 0.9692665404913205
$ws.Range("E3").Value = 0.9845133521142924
$ws.Range("F3").Value = 1.006200714633619
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = -0.03231471279853988
$ws.Range("C4").Value = 0.7695157101427289
$ws.Range("D4").Value = 0.9504010021523444
$ws.Range("E4").Value = 0.9748851225412891
$ws.Range("F4").Value = 1.004336377456358
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.05671622629984657
$ws.Range("C5").Value = 0.6321910596607474
$ws.Range("D5").Value = 0.8518300650585396
$ws.Range("E5").Value = 0.9229464042177854
$ws.Range("F5").Value = 0.951413452211165
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.085928964333323
$ws.Range("C6").Value = 0.5612671955956002
$ws.Range("D6").Value = 0.5873845025919633
$ws.Range("E6").Value = 0.7664101399328973
$ws.Range("F6").Value = 0.7883078955770884
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.1382138394267656
$ws.Range("C7").Value = 0.7476276879241753
$ws.Range("D7").Value = 1.097975515794379
$ws.Range("E7").Value = 1.047843268716452
$ws.Range("F7").Value = 1.077897322974905
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.1215206328360093
$ws.Range("C8").Value = 0.7819526689445239
$ws.Range("D8").Value = 1.205421270947194
$ws.Range("E8").Value = 1.097916786895616
$ws.Range("F8").Value = 1.135726716235405
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.2028461735751207
$ws.Range("C9").Value = 0.8472735835465873
$ws.Range("D9").Value = 1.158512388347609
$ws.Range("E9").Value = 1.076342133500128
$ws.Range("F9").Value = 1.10405820904526
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.239907996146195
$ws.Range("C10").Value = 0.6651045147971374
$ws.Range("D10").Value = 0.5928596743010739
$ws.Range("E10").Value = 0.769973814035954
$ws.Range("F10").Value = 0.7673553351966808
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.2225839694609408
$ws.Range("C11").Value = 0.8306889768957216
$ws.Range("D11").Value = 1.245987321698922
$ws.Range("E11").Value = 1.116238021973325
$ws.Range("F11").Value = 1.152988242321246
$ws.Range("G11").Value = 10

